$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'22.993.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.78%  "

$ws.Range("D3").Value = "'1.582.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.86%  "

$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.39%  "

$ws.Range("E5").Value = "  +0.51%  "

$ws.Range("D6").Value = "'298.92"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.22%  "

$ws.Range("D7").Value = "'0.3746"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.90%  "

$ws.Range("D8").Value = "'0.3556"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.30%  "

$ws.Range("D9").Value = "'50.48"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.70%  "

$ws.Range("D10").Value = "'1.004"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.35%  "

$ws.Range("D11").Value = "'1.209"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.59%  "

$ws.Range("D12").Value = "'0.07953"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.95%  "

$ws.Range("D13").Value = "'21.70"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.36%  "

$ws.Range("D14").Value = "'6.432"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.24%  "

$ws.Range("D15").Value = "'7.253"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.35%  "

$ws.Range("D16").Value = "'0.00001213"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.86%  "

$ws.Range("D17").Value = "'1.586.28"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.53%  "

$ws.Range("D18").Value = "'91.71"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.05%  "

$ws.Range("D19").Value = "'0.06741"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.86%  "

$ws.Range("D20").Value = "'17.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.26%  "

$ws.Range("D21").Value = "'1.003"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.52%  "

$ws.Range("D22").Value = "'6.371"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.60%  "

$ws.Range("B23").Value = "WrappedBTC"
$ws.Range("C23").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D23").Value = "'22.997.72"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.82%  "

$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").Value = "'12.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.36%  "

$ws.Range("D25").Value = "'2.370"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.34%  "

$ws.Range("D26").Value = "'2.757"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.84%  "

$ws.Range("D27").Value = "'20.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.92%  "

$ws.Range("D28").Value = "'146.51"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.78%  "

$ws.Range("D29").Value = "'5.204"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.40%  "

$ws.Range("D30").Value = "'131.28"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.21%  "

$ws.Range("E31").Value = "  -4.81%  "

$ws.Range("B32").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C32").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D32").Value = "'1.767.50"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.15%  "

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'6.441"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.73%  "

$ws.Range("D34").Value = "'0.9248"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.84%  "

$ws.Range("D35").Value = "'0.07296"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.02%  "

$ws.Range("D36").Value = "'0.02651"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.20%  "

$ws.Range("D37").Value = "'0.08727"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.72%  "

$ws.Range("E38").Value = "  -3.24%  "

$ws.Range("D39").Value = "'9.815"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.30%  "

$ws.Range("D40").Value = "'5.963"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.02%  "

$ws.Range("E41").Value = "  -5.18%  "

$ws.Range("D42").Value = "'0.6812"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.53%  "

$ws.Range("D43").Value = "'11.65"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -9.36%  "

$ws.Range("D44").Value = "'14.71"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.53%  "

$ws.Range("D45").Value = "'0.6301"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.19%  "

$ws.Range("D46").Value = "'3.962"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.42%  "

$ws.Range("D47").Value = "'2.225"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.98%  "

$ws.Range("D48").Value = "'130.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.22%  "

$ws.Range("D49").Value = "'0.07865"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.09%  "

$ws.Range("D50").Value = "'1.174"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.01%  "

$ws.Range("D51").Value = "'1.160"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.59%  "
